$wb = $excel.ActiveWorkbook

# Sheet ALC (sheet1.xml)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 289297.78  # H17: 306731.22 -> 289297.78
$ws.Cells.Item(17, 10).Value = 289297.78  # J17: 306731.22 -> 289297.78
$ws.Cells.Item(17, 12).Value = 867893.3400000001  # L17: 920193.6599999999 -> 867893.3400000001
$ws.Cells.Item(17, 14).Value = -868229.3400000001  # N17: -920529.6599999999 -> -868229.3400000001
$ws.Cells.Item(86, 8).Value = 4626.9473  # H86: 4169.696 -> 4626.9473
$ws.Cells.Item(86, 10).Value = 6101.091  # J86: 5006.8667 -> 6101.091
$ws.Cells.Item(86, 12).Value = 6101.091  # L86: 5006.8667 -> 6101.091
$ws.Cells.Item(86, 14).Value = -8347.091  # N86: -7252.8667 -> -8347.091
$ws.Cells.Item(89, 8).Value = 4626.9473  # H89: 4169.696 -> 4626.9473
$ws.Cells.Item(89, 10).Value = 6101.091  # J89: 5006.8667 -> 6101.091
$ws.Cells.Item(89, 12).Value = 30505.455  # L89: 25034.3335 -> 30505.455
$ws.Cells.Item(89, 14).Value = -41737.455  # N89: -36266.33349999999 -> -41737.455
$ws.Cells.Item(92, 8).Value = 2605.7  # H92: 2690.2632 -> 2605.7
$ws.Cells.Item(92, 9).Value = 2105.5386  # I92: 2197.75 -> 2105.5386
$ws.Cells.Item(92, 11).Value = 2105.5386  # K92: 2197.75 -> 2105.5386
$ws.Cells.Item(92, 13).Value = -857.5385999999999  # M92: -949.75 -> -857.5385999999999
$ws.Cells.Item(106, 8).Value = 3698.4443  # H106: 3548.6 -> 3698.4443
$ws.Cells.Item(106, 9).Value = 3160.75  # I106: 3054 -> 3160.75
$ws.Cells.Item(106, 11).Value = 3160.75  # K106: 3054 -> 3160.75
$ws.Cells.Item(106, 13).Value = -2529.75  # M106: -2423 -> -2529.75
$ws.Cells.Item(111, 8).Value = 16580.63  # H111: 17930.111 -> 16580.63
$ws.Cells.Item(111, 9).Value = 17131.295  # I111: 18683.875 -> 17131.295
$ws.Cells.Item(111, 11).Value = 51393.88499999999  # K111: 56051.625 -> 51393.88499999999
$ws.Cells.Item(111, 13).Value = -48326.88499999999  # M111: -52984.625 -> -48326.88499999999
$ws.Cells.Item(112, 8).Value = 836349.8  # H112: 772169.0600000001 -> 836349.8
$ws.Cells.Item(112, 10).Value = 1003249.9  # J112: 912227.2 -> 1003249.9
$ws.Cells.Item(112, 12).Value = 3009749.7  # L112: 2736681.6 -> 3009749.7
$ws.Cells.Item(112, 14).Value = -3011965.7  # N112: -2738897.6 -> -3011965.7
$ws.Cells.Item(127, 8).Value = 845.125  # H127: 923.2857 -> 845.125
$ws.Cells.Item(127, 9).Value = 845.125  # I127: 923.2857 -> 845.125
$ws.Cells.Item(127, 11).Value = 2535.375  # K127: 2769.8571 -> 2535.375
$ws.Cells.Item(127, 13).Value = 2424.625  # M127: 2190.1429 -> 2424.625
$ws.Cells.Item(129, 8).Value = 100006320  # H129: 111118024 -> 100006320
$ws.Cells.Item(129, 9).Value = 200000800  # I129: 250000750 -> 200000800
$ws.Cells.Item(129, 11).Value = 600002400  # K129: 750002250 -> 600002400
$ws.Cells.Item(129, 13).Value = -599997400  # M129: -749997250 -> -599997400
$ws.Cells.Item(132, 8).Value = 43482556  # H132: 47624504 -> 43482556
$ws.Cells.Item(132, 9).Value = 50004716  # I132: 55561670 -> 50004716
$ws.Cells.Item(132, 11).Value = 150014148  # K132: 166685010 -> 150014148
$ws.Cells.Item(132, 13).Value = -150011618  # M132: -166682480 -> -150011618
$ws.Cells.Item(137, 8).Value = 4946.067  # H137: 5282.2925 -> 4946.067
$ws.Cells.Item(137, 9).Value = 2847.3704  # I137: 3015.7917 -> 2847.3704
$ws.Cells.Item(137, 10).Value = 8094.1113  # J137: 8482.058999999999 -> 8094.1113
$ws.Cells.Item(137, 11).Value = 8542.111199999999  # K137: 9047.375100000001 -> 8542.111199999999
$ws.Cells.Item(137, 12).Value = 24282.3339  # L137: 25446.177 -> 24282.3339
$ws.Cells.Item(137, 13).Value = -5992.111199999999  # M137: -6497.375100000001 -> -5992.111199999999
$ws.Cells.Item(137, 14).Value = -29382.3339  # N137: -30546.177 -> -29382.3339
$ws.Cells.Item(138, 8).Value = 113408.71  # H138: 137370.06 -> 113408.71
$ws.Cells.Item(138, 9).Value = 1262.6111  # I138: 0 -> 1262.6111
$ws.Cells.Item(138, 10).Value = 141061.17  # J138: 137370.06 -> 141061.17
$ws.Cells.Item(138, 11).Value = 3787.8333  # K138: 0 -> 3787.8333
$ws.Cells.Item(138, 12).Value = 423183.51  # L138: 412110.18 -> 423183.51
$ws.Cells.Item(138, 13).Value = 1352.1667  # M138: None -> 1352.1667
$ws.Cells.Item(138, 14).Value = -433463.51  # N138: -422390.18 -> -433463.51
$ws.Cells.Item(141, 8).Value = 4662.52  # H141: 4769.5 -> 4662.52
$ws.Cells.Item(141, 9).Value = 2611.4546  # I141: 2663.1 -> 2611.4546
$ws.Cells.Item(141, 11).Value = 7834.3638  # K141: 7989.299999999999 -> 7834.3638
$ws.Cells.Item(141, 13).Value = -2654.3638  # M141: -2809.299999999999 -> -2654.3638

# Sheet ARM (sheet2.xml)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8609.013999999999  # H32: 8988.973 -> 8609.013999999999
$ws.Cells.Item(32, 9).Value = 7895.471  # I32: 8271.835999999999 -> 7895.471
$ws.Cells.Item(32, 11).Value = 7895.471  # K32: 8271.835999999999 -> 7895.471
$ws.Cells.Item(32, 13).Value = -7608.471  # M32: -7984.835999999999 -> -7608.471
$ws.Cells.Item(61, 8).Value = 6130.88  # H61: 7244.4834 -> 6130.88
$ws.Cells.Item(61, 9).Value = 4301.875  # I61: 4852.321 -> 4301.875
$ws.Cells.Item(61, 10).Value = 16772.363  # J61: 25356.572 -> 16772.363
$ws.Cells.Item(61, 11).Value = 4301.875  # K61: 4852.321 -> 4301.875
$ws.Cells.Item(61, 12).Value = 16772.363  # L61: 25356.572 -> 16772.363
$ws.Cells.Item(61, 13).Value = -4089.875  # M61: -4640.321 -> -4089.875
$ws.Cells.Item(61, 14).Value = -17196.363  # N61: -25780.572 -> -17196.363
$ws.Cells.Item(74, 8).Value = 2335.3262  # H74: 2434.0466 -> 2335.3262
$ws.Cells.Item(74, 9).Value = 1407.1538  # I74: 1447.7222 -> 1407.1538
$ws.Cells.Item(74, 11).Value = 1407.1538  # K74: 1447.7222 -> 1407.1538
$ws.Cells.Item(74, 13).Value = -533.1538  # M74: -573.7221999999999 -> -533.1538
$ws.Cells.Item(77, 8).Value = 2335.3262  # H77: 2434.0466 -> 2335.3262
$ws.Cells.Item(77, 9).Value = 1407.1538  # I77: 1447.7222 -> 1407.1538
$ws.Cells.Item(77, 11).Value = 7035.769  # K77: 7238.611 -> 7035.769
$ws.Cells.Item(77, 13).Value = -2667.769  # M77: -2870.611 -> -2667.769
$ws.Cells.Item(97, 8).Value = 737.37836  # H97: 748.1389 -> 737.37836
$ws.Cells.Item(97, 9).Value = 683.5  # I97: 694.25806 -> 683.5
$ws.Cells.Item(97, 11).Value = 683.5  # K97: 694.25806 -> 683.5
$ws.Cells.Item(97, 13).Value = -187.5  # M97: -198.25806 -> -187.5
$ws.Cells.Item(110, 8).Value = 888.41174  # H110: 1048.15 -> 888.41174
$ws.Cells.Item(110, 9).Value = 663.7143  # I110: 833.6875 -> 663.7143
$ws.Cells.Item(110, 10).Value = 1937  # J110: 1906 -> 1937
$ws.Cells.Item(110, 11).Value = 663.7143  # K110: 833.6875 -> 663.7143
$ws.Cells.Item(110, 12).Value = 1937  # L110: 1906 -> 1937
$ws.Cells.Item(110, 13).Value = 1381.2857  # M110: 1211.3125 -> 1381.2857
$ws.Cells.Item(110, 14).Value = -6027  # N110: -5996 -> -6027
$ws.Cells.Item(123, 8).Value = 80475  # H123: 82025 -> 80475
$ws.Cells.Item(123, 9).Value = 0  # I123: 77100 -> 0
$ws.Cells.Item(123, 10).Value = 80475  # J123: 83666.664 -> 80475
$ws.Cells.Item(123, 11).Value = 0  # K123: 77100 -> 0
$ws.Cells.Item(123, 12).Value = 80475  # L123: 83666.664 -> 80475
$ws.Cells.Item(123, 13).ClearContents()  # M123: delete (was -72200)
$ws.Cells.Item(123, 14).Value = -90275  # N123: -93466.664 -> -90275
$ws.Cells.Item(132, 8).Value = 1979.125  # H132: 2216.205 -> 1979.125
$ws.Cells.Item(132, 9).Value = 1923.2307  # I132: 2175.484 -> 1923.2307
$ws.Cells.Item(132, 10).Value = 2221.3333  # J132: 2374 -> 2221.3333
$ws.Cells.Item(132, 11).Value = 5769.6921  # K132: 6526.451999999999 -> 5769.6921
$ws.Cells.Item(132, 12).Value = 6663.999899999999  # L132: 7122 -> 6663.999899999999
$ws.Cells.Item(132, 13).Value = -3239.6921  # M132: -3996.451999999999 -> -3239.6921
$ws.Cells.Item(132, 14).Value = -11723.9999  # N132: -12182 -> -11723.9999
$ws.Cells.Item(136, 8).Value = 6130.88  # H136: 7244.4834 -> 6130.88
$ws.Cells.Item(136, 9).Value = 4301.875  # I136: 4852.321 -> 4301.875
$ws.Cells.Item(136, 10).Value = 16772.363  # J136: 25356.572 -> 16772.363
$ws.Cells.Item(136, 11).Value = 12905.625  # K136: 14556.963 -> 12905.625
$ws.Cells.Item(136, 12).Value = 50317.08900000001  # L136: 76069.716 -> 50317.08900000001
$ws.Cells.Item(136, 13).Value = -10355.625  # M136: -12006.963 -> -10355.625
$ws.Cells.Item(136, 14).Value = -55417.08900000001  # N136: -81169.716 -> -55417.08900000001

# Sheet BSM (sheet3.xml)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2609.513  # H20: 2633.6052 -> 2609.513
$ws.Cells.Item(20, 9).Value = 2509.1936  # I20: 2536.3667 -> 2509.1936
$ws.Cells.Item(20, 11).Value = 2509.1936  # K20: 2536.3667 -> 2509.1936
$ws.Cells.Item(20, 13).Value = -2262.1936  # M20: -2289.3667 -> -2262.1936
$ws.Cells.Item(94, 8).Value = 1274.9375  # H94: 1296.875 -> 1274.9375
$ws.Cells.Item(94, 9).Value = 1407.4615  # I94: 1396.0741 -> 1407.4615
$ws.Cells.Item(94, 10).Value = 700.6667  # J94: 761.2 -> 700.6667
$ws.Cells.Item(94, 11).Value = 1407.4615  # K94: 1396.0741 -> 1407.4615
$ws.Cells.Item(94, 12).Value = 700.6667  # L94: 761.2 -> 700.6667
$ws.Cells.Item(94, 13).Value = -956.4614999999999  # M94: -945.0741 -> -956.4614999999999
$ws.Cells.Item(94, 14).Value = -1602.6667  # N94: -1663.2 -> -1602.6667
$ws.Cells.Item(105, 8).Value = 4861.409  # H105: 5195.1 -> 4861.409
$ws.Cells.Item(105, 9).Value = 3812.5  # I105: 4193.8335 -> 3812.5
$ws.Cells.Item(105, 11).Value = 3812.5  # K105: 4193.8335 -> 3812.5
$ws.Cells.Item(105, 13).Value = -2065.5  # M105: -2446.8335 -> -2065.5
$ws.Cells.Item(107, 8).Value = 1663.9048  # H107: 1664.381 -> 1663.9048
$ws.Cells.Item(107, 9).Value = 1418.3572  # I107: 1419.0714 -> 1418.3572
$ws.Cells.Item(107, 11).Value = 1418.3572  # K107: 1419.0714 -> 1418.3572
$ws.Cells.Item(107, 13).Value = 501.6428000000001  # M107: 500.9286 -> 501.6428000000001
$ws.Cells.Item(134, 8).Value = 2201.743  # H134: 2355 -> 2201.743
$ws.Cells.Item(134, 9).Value = 2073.9673  # I134: 2170.9312 -> 2073.9673
$ws.Cells.Item(134, 10).Value = 3067.7778  # J134: 3880.1428 -> 3067.7778
$ws.Cells.Item(134, 11).Value = 6221.901899999999  # K134: 6512.7936 -> 6221.901899999999
$ws.Cells.Item(134, 12).Value = 9203.3334  # L134: 11640.4284 -> 9203.3334
$ws.Cells.Item(134, 13).Value = -3686.901899999999  # M134: -3977.7936 -> -3686.901899999999
$ws.Cells.Item(134, 14).Value = -14273.3334  # N134: -16710.4284 -> -14273.3334

# Sheet CRP (sheet4.xml)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1608.8518  # H16: 1744.6 -> 1608.8518
$ws.Cells.Item(16, 9).Value = 1056.1578  # I16: 1230.7778 -> 1056.1578
$ws.Cells.Item(16, 10).Value = 2921.5  # J16: 3065.8572 -> 2921.5
$ws.Cells.Item(16, 11).Value = 1056.1578  # K16: 1230.7778 -> 1056.1578
$ws.Cells.Item(16, 12).Value = 2921.5  # L16: 3065.8572 -> 2921.5
$ws.Cells.Item(16, 13).Value = -769.1578  # M16: -943.7778000000001 -> -769.1578
$ws.Cells.Item(16, 14).Value = -3495.5  # N16: -3639.8572 -> -3495.5
$ws.Cells.Item(44, 8).Value = 29500  # H44: 28000 -> 29500
$ws.Cells.Item(44, 10).Value = 30000  # J44: 0 -> 30000
$ws.Cells.Item(44, 12).Value = 30000  # L44: 0 -> 30000
$ws.Cells.Item(44, 14).Value = -30884  # N44: None -> -30884
$ws.Cells.Item(113, 8).Value = 1608.8518  # H113: 1744.6 -> 1608.8518
$ws.Cells.Item(113, 9).Value = 1056.1578  # I113: 1230.7778 -> 1056.1578
$ws.Cells.Item(113, 10).Value = 2921.5  # J113: 3065.8572 -> 2921.5
$ws.Cells.Item(113, 11).Value = 1056.1578  # K113: 1230.7778 -> 1056.1578
$ws.Cells.Item(113, 12).Value = 2921.5  # L113: 3065.8572 -> 2921.5
$ws.Cells.Item(113, 13).Value = 1113.8422  # M113: 939.2221999999999 -> 1113.8422
$ws.Cells.Item(113, 14).Value = -7261.5  # N113: -7405.8572 -> -7261.5
$ws.Cells.Item(132, 8).Value = 2108331.2  # H132: 2225361 -> 2108331.2
$ws.Cells.Item(132, 9).Value = 3079099.8  # I132: 3335541.2 -> 3079099.8
$ws.Cells.Item(132, 11).Value = 9237299.399999999  # K132: 10006623.6 -> 9237299.399999999
$ws.Cells.Item(132, 13).Value = -9234769.399999999  # M132: -10004093.6 -> -9234769.399999999
$ws.Cells.Item(134, 8).Value = 4492.3  # H134: 4430.2 -> 4492.3
$ws.Cells.Item(134, 9).Value = 2070.5  # I134: 1967 -> 2070.5
$ws.Cells.Item(134, 11).Value = 6211.5  # K134: 5901 -> 6211.5
$ws.Cells.Item(134, 13).Value = -3676.5  # M134: -3366 -> -3676.5
$ws.Cells.Item(141, 8).Value = 678040.5600000001  # H141: 567027.4399999999 -> 678040.5600000001
$ws.Cells.Item(141, 10).Value = 678040.5600000001  # J141: 567027.4399999999 -> 678040.5600000001
$ws.Cells.Item(141, 12).Value = 678040.5600000001  # L141: 567027.4399999999 -> 678040.5600000001
$ws.Cells.Item(141, 14).Value = -688400.5600000001  # N141: -577387.4399999999 -> -688400.5600000001

# Sheet CUL (sheet5.xml)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 117935.9  # H131: 162904.5 -> 117935.9
$ws.Cells.Item(131, 10).Value = 1744  # J131: 1733.585 -> 1744
$ws.Cells.Item(131, 12).Value = 5232  # L131: 5200.755 -> 5232
$ws.Cells.Item(131, 14).Value = -15312  # N131: -15280.755 -> -15312

# Sheet GSM (sheet6.xml)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 87000  # H39: 90000 -> 87000
$ws.Cells.Item(39, 10).Value = 87000  # J39: 90000 -> 87000
$ws.Cells.Item(39, 12).Value = 87000  # L39: 90000 -> 87000
$ws.Cells.Item(39, 14).Value = -88064  # N39: -91064 -> -88064
$ws.Cells.Item(58, 8).Value = 0  # H58: 35555 -> 0
$ws.Cells.Item(58, 10).Value = 0  # J58: 35555 -> 0
$ws.Cells.Item(58, 12).Value = 0  # L58: 35555 -> 0
$ws.Cells.Item(58, 14).ClearContents()  # N58: delete (was -36109)
$ws.Cells.Item(80, 8).Value = 4032.2856  # H80: 3959.3125 -> 4032.2856
$ws.Cells.Item(80, 9).Value = 3772  # I80: 3710.111 -> 3772
$ws.Cells.Item(80, 10).Value = 4292.5713  # J80: 4279.7144 -> 4292.5713
$ws.Cells.Item(80, 11).Value = 3772  # K80: 3710.111 -> 3772
$ws.Cells.Item(80, 12).Value = 4292.5713  # L80: 4279.7144 -> 4292.5713
$ws.Cells.Item(80, 13).Value = -2774  # M80: -2712.111 -> -2774
$ws.Cells.Item(80, 14).Value = -6288.5713  # N80: -6275.7144 -> -6288.5713
$ws.Cells.Item(83, 8).Value = 4032.2856  # H83: 3959.3125 -> 4032.2856
$ws.Cells.Item(83, 9).Value = 3772  # I83: 3710.111 -> 3772
$ws.Cells.Item(83, 10).Value = 4292.5713  # J83: 4279.7144 -> 4292.5713
$ws.Cells.Item(83, 11).Value = 18860  # K83: 18550.555 -> 18860
$ws.Cells.Item(83, 12).Value = 21462.8565  # L83: 21398.572 -> 21462.8565
$ws.Cells.Item(83, 13).Value = -13868  # M83: -13558.555 -> -13868
$ws.Cells.Item(83, 14).Value = -31446.8565  # N83: -31382.572 -> -31446.8565
$ws.Cells.Item(97, 8).Value = 949.9032  # H97: 954.625 -> 949.9032
$ws.Cells.Item(97, 9).Value = 613.86957  # I97: 624.36 -> 613.86957
$ws.Cells.Item(97, 10).Value = 1916  # J97: 2134.1428 -> 1916
$ws.Cells.Item(97, 11).Value = 613.86957  # K97: 624.36 -> 613.86957
$ws.Cells.Item(97, 12).Value = 1916  # L97: 2134.1428 -> 1916
$ws.Cells.Item(97, 13).Value = -117.86957  # M97: -128.36 -> -117.86957
$ws.Cells.Item(97, 14).Value = -2908  # N97: -3126.1428 -> -2908
$ws.Cells.Item(132, 8).Value = 8132888  # H132: 8336178 -> 8132888
$ws.Cells.Item(132, 10).Value = 2844.5  # J132: 3363 -> 2844.5
$ws.Cells.Item(132, 12).Value = 8533.5  # L132: 10089 -> 8533.5
$ws.Cells.Item(132, 14).Value = -13593.5  # N132: -15149 -> -13593.5
$ws.Cells.Item(136, 8).Value = 73660  # H136: 72161 -> 73660
$ws.Cells.Item(136, 10).Value = 73660  # J136: 72161 -> 73660
$ws.Cells.Item(136, 12).Value = 220980  # L136: 216483 -> 220980
$ws.Cells.Item(136, 14).Value = -226080  # N136: -221583 -> -226080

# Sheet LTW (sheet7.xml)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7024.7334  # H7: 6519.514 -> 7024.7334
$ws.Cells.Item(7, 9).Value = 6275.095  # I7: 6383.1577 -> 6275.095
$ws.Cells.Item(7, 10).Value = 8773.888999999999  # J7: 6681.4375 -> 8773.888999999999
$ws.Cells.Item(7, 11).Value = 6275.095  # K7: 6383.1577 -> 6275.095
$ws.Cells.Item(7, 12).Value = 8773.888999999999  # L7: 6681.4375 -> 8773.888999999999
$ws.Cells.Item(7, 13).Value = -6163.095  # M7: -6271.1577 -> -6163.095
$ws.Cells.Item(7, 14).Value = -8997.888999999999  # N7: -6905.4375 -> -8997.888999999999
$ws.Cells.Item(68, 8).Value = 2830.6924  # H68: 2830.7693 -> 2830.6924
$ws.Cells.Item(68, 9).Value = 2690.3  # I68: 2690.4 -> 2690.3
$ws.Cells.Item(68, 11).Value = 2690.3  # K68: 2690.4 -> 2690.3
$ws.Cells.Item(68, 13).Value = -1941.3  # M68: -1941.4 -> -1941.3
$ws.Cells.Item(71, 8).Value = 2830.6924  # H71: 2830.7693 -> 2830.6924
$ws.Cells.Item(71, 9).Value = 2690.3  # I71: 2690.4 -> 2690.3
$ws.Cells.Item(71, 11).Value = 13451.5  # K71: 13452 -> 13451.5
$ws.Cells.Item(71, 13).Value = -9707.5  # M71: -9708 -> -9707.5
$ws.Cells.Item(93, 8).Value = 3268.64  # H93: 3686 -> 3268.64
$ws.Cells.Item(93, 9).Value = 3720.5557  # I93: 4159.375 -> 3720.5557
$ws.Cells.Item(93, 10).Value = 2106.5715  # J93: 2423.6667 -> 2106.5715
$ws.Cells.Item(93, 11).Value = 3720.5557  # K93: 4159.375 -> 3720.5557
$ws.Cells.Item(93, 12).Value = 2106.5715  # L93: 2423.6667 -> 2106.5715
$ws.Cells.Item(93, 13).Value = -2472.5557  # M93: -2911.375 -> -2472.5557
$ws.Cells.Item(93, 14).Value = -4602.5715  # N93: -4919.6667 -> -4602.5715
$ws.Cells.Item(126, 8).Value = 7024.7334  # H126: 6519.514 -> 7024.7334
$ws.Cells.Item(126, 9).Value = 6275.095  # I126: 6383.1577 -> 6275.095
$ws.Cells.Item(126, 10).Value = 8773.888999999999  # J126: 6681.4375 -> 8773.888999999999
$ws.Cells.Item(126, 11).Value = 18825.285  # K126: 19149.4731 -> 18825.285
$ws.Cells.Item(126, 12).Value = 26321.667  # L126: 20044.3125 -> 26321.667
$ws.Cells.Item(126, 13).Value = -16355.285  # M126: -16679.4731 -> -16355.285
$ws.Cells.Item(126, 14).Value = -31261.667  # N126: -24984.3125 -> -31261.667

# Sheet WVR (sheet8.xml)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 29998  # H24: 0 -> 29998
$ws.Cells.Item(24, 10).Value = 29998  # J24: 0 -> 29998
$ws.Cells.Item(24, 12).Value = 29998  # L24: 0 -> 29998
$ws.Cells.Item(24, 14).Value = -30458  # N24: None -> -30458
$ws.Cells.Item(122, 8).Value = 1506.72  # H122: 1606.6957 -> 1506.72
$ws.Cells.Item(122, 9).Value = 1606.6522  # I122: 1663.3636 -> 1606.6522
$ws.Cells.Item(122, 10).Value = 357.5  # J122: 360 -> 357.5
$ws.Cells.Item(122, 11).Value = 4819.9566  # K122: 4990.0908 -> 4819.9566
$ws.Cells.Item(122, 12).Value = 1072.5  # L122: 1080 -> 1072.5
$ws.Cells.Item(122, 13).Value = -2369.9566  # M122: -2540.0908 -> -2369.9566
$ws.Cells.Item(122, 14).Value = -5972.5  # N122: -5980 -> -5972.5
$ws.Cells.Item(136, 8).Value = 6839.951  # H136: 6994.1 -> 6839.951
$ws.Cells.Item(136, 9).Value = 7239.2354  # I136: 7438.1816 -> 7239.2354
$ws.Cells.Item(136, 11).Value = 21717.7062  # K136: 22314.5448 -> 21717.7062
$ws.Cells.Item(136, 13).Value = -19167.7062  # M136: -19764.5448 -> -19167.7062
